$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column A (ID, AGE), shifting the existing
# header/data columns right by two (A->C, B->D, ... H->J).
$ws.Range("A:B").Insert()

# New header cells, matching the bold/centered style used by the other headers
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "AGE"
$ws.Range("C1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2 (existing person) - fill in auto-incrementing ID and AGE
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 18

# New rows for additional people, each with an auto-incrementing ID
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 22
$ws.Range("F3").Value = 44.78070068359375

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 25
$ws.Range("F4").Value = 46.16134643554688

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 33
$ws.Range("F5").Value = 40.92935180664062
